$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the label text in D1 (identity/status column header)
$ws.Range("D1").Value = "身分 (學士、碩士或博士班）"

# Widen column D to fit the new, longer text (target stored width ~30.5)
$ws.Columns.Item(4).ColumnWidth = 29.8

# Move the active selection to G5, matching the saved view state
$ws.Range("G5").Select()
